# Fixed weighting for ccpi data
# Columns: H = 2024_CCPI_GHG, I = 2024_CCPI_Renewables,
#          J = 2024_CCPI_Energy_Use, K = 2024_CCPI_Climate_Policy
#
# The CCPI sub-indicator columns were being under-weighted; rescale them
# to the corrected weights (GHG x2.5, Renewables/Energy Use/Climate Policy x5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$ghgFactor = 2.5
$otherFactor = 5

for ($r = 2; $r -le $lastRow; $r++) {
    $hCell = $ws.Cells.Item($r, 8)
    $hVal = $hCell.Value2

    if ($hVal -ne $null -and $hVal -ne "") {
        $iCell = $ws.Cells.Item($r, 9)
        $jCell = $ws.Cells.Item($r, 10)
        $kCell = $ws.Cells.Item($r, 11)

        $hCell.Value2 = $hVal * $ghgFactor
        $iCell.Value2 = $iCell.Value2 * $otherFactor
        $jCell.Value2 = $jCell.Value2 * $otherFactor
        $kCell.Value2 = $kCell.Value2 * $otherFactor
    }
}
